$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "D" column holds dates stored as literal text (e.g. "01/01/2023").
# When such a text is assigned directly to .Value, Excel auto-converts it
# into a real date serial number. To preserve the original text semantics
# we force the cell to Text format before assigning, then restore the
# cell style back to "Normal" so no stray formatting is left behind.

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 32
Set-TextValue $ws.Range("C32") "Administração Pública"
Set-TextValue $ws.Range("D32") "01/01/2022"
$ws.Range("E32").Value = 2.409363299680759

# Row 33
Set-TextValue $ws.Range("C33") "Entidades Empresariais"
Set-TextValue $ws.Range("D33") "01/01/2022"
$ws.Range("E33").Value = 40.19065623217473

# Row 34
Set-TextValue $ws.Range("C34") "Entidades sem Fins Lucrativos"
Set-TextValue $ws.Range("D34") "01/01/2022"
$ws.Range("E34").Value = 22.19863338455868

# Row 35
Set-TextValue $ws.Range("C35") "Administração Pública"
Set-TextValue $ws.Range("D35") "01/01/2023"
$ws.Range("E35").Value = 2.315008864560557

# Row 36
Set-TextValue $ws.Range("C36") "Entidades Empresariais"
Set-TextValue $ws.Range("D36") "01/01/2023"
$ws.Range("E36").Value = 40.48422032949609
